$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 999
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H76").Value = 4997.6665
$ws.Range("I76").Value = 4997.6665
$ws.Range("K76").Value = 4997.6665
$ws.Range("M76").Value = -4682.6665
$ws.Range("H79").Value = 4997.6665
$ws.Range("I79").Value = 4997.6665
$ws.Range("K79").Value = 4997.6665
$ws.Range("M79").Value = -3905.6665
$ws.Range("H135").Value = 1380.909
$ws.Range("I135").Value = 1036.375
$ws.Range("J135").Value = 2299.6667
$ws.Range("K135").Value = 9327.375
$ws.Range("L135").Value = 20697.0003
$ws.Range("M135").Value = -6792.375
$ws.Range("N135").Value = -25767.0003
$ws.Range("H137").Value = 1583.1333
$ws.Range("I137").Value = 1431.6364
$ws.Range("K137").Value = 4294.9092
$ws.Range("M137").Value = -1744.9092

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1823.6428
$ws.Range("I2").Value = 593.8
$ws.Range("J2").Value = 4898.25
$ws.Range("K2").Value = 593.8
$ws.Range("L2").Value = 4898.25
$ws.Range("M2").Value = -480.8
$ws.Range("N2").Value = -5124.25
$ws.Range("H44").Value = 47993.25
$ws.Range("J44").Value = 47993.25
$ws.Range("L44").Value = 47993.25
$ws.Range("N44").Value = -48969.25
$ws.Range("H55").Value = 26670.166
$ws.Range("I55").Value = 12024
$ws.Range("J55").Value = 33993.25
$ws.Range("K55").Value = 12024
$ws.Range("L55").Value = 33993.25
$ws.Range("M55").Value = -11709
$ws.Range("N55").Value = -34623.25
$ws.Range("H80").Value = 89947.5
$ws.Range("J80").Value = 89947.5
$ws.Range("L80").Value = 89947.5
$ws.Range("N80").Value = -91943.5
$ws.Range("H83").Value = 89947.5
$ws.Range("J83").Value = 89947.5
$ws.Range("L83").Value = 269842.5
$ws.Range("N83").Value = -279826.5
$ws.Range("H116").Value = 1823.6428
$ws.Range("I116").Value = 593.8
$ws.Range("J116").Value = 4898.25
$ws.Range("K116").Value = 593.8
$ws.Range("L116").Value = 4898.25
$ws.Range("M116").Value = 1700.2
$ws.Range("N116").Value = -9486.25
$ws.Range("H122").Value = 3023.625
$ws.Range("I122").Value = 2138.2
$ws.Range("K122").Value = 6414.599999999999
$ws.Range("M122").Value = -3964.599999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1823.6428
$ws.Range("I3").Value = 593.8
$ws.Range("J3").Value = 4898.25
$ws.Range("K3").Value = 593.8
$ws.Range("L3").Value = 4898.25
$ws.Range("M3").Value = -479.8
$ws.Range("N3").Value = -5126.25
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 2000
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -4246
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 10000
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -21232

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3133.3333
$ws.Range("J16").Value = 3999
$ws.Range("L16").Value = 3999
$ws.Range("N16").Value = -4573
$ws.Range("H58").Value = 1963.8125
$ws.Range("I58").Value = 1100.76
$ws.Range("K58").Value = 1100.76
$ws.Range("M58").Value = -897.76
$ws.Range("H113").Value = 3133.3333
$ws.Range("J113").Value = 3999
$ws.Range("L113").Value = 3999
$ws.Range("N113").Value = -8339
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 1877.75
$ws.Range("I132").Value = 1877.75
$ws.Range("K132").Value = 5633.25
$ws.Range("M132").Value = -3103.25
$ws.Range("H136").Value = 1963.8125
$ws.Range("I136").Value = 1100.76
$ws.Range("K136").Value = 3302.28
$ws.Range("M136").Value = -752.2799999999997

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2400
$ws.Range("J39").Value = 2400
$ws.Range("L39").Value = 7200
$ws.Range("N39").Value = -7788
$ws.Range("H59").Value = 1005
$ws.Range("I59").Value = 1005
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 3015
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -2475
$ws.Range("N59").ClearContents()
$ws.Range("H96").Value = 15000
$ws.Range("J96").Value = 15000
$ws.Range("L96").Value = 45000
$ws.Range("N96").Value = -49118
$ws.Range("H97").Value = 8693
$ws.Range("I97").Value = 7995
$ws.Range("J97").Value = 9158.333000000001
$ws.Range("K97").Value = 23985
$ws.Range("L97").Value = 27474.999
$ws.Range("M97").Value = -23489
$ws.Range("N97").Value = -28466.999
$ws.Range("H129").Value = 1397.75
$ws.Range("J129").Value = 2500
$ws.Range("L129").Value = 7500
$ws.Range("N129").Value = -17500
$ws.Range("H130").Value = 2997.1428
$ws.Range("I130").Value = 1666
$ws.Range("J130").Value = 3995.5
$ws.Range("K130").Value = 4998
$ws.Range("L130").Value = 11986.5
$ws.Range("M130").Value = 22
$ws.Range("N130").Value = -22026.5
$ws.Range("H131").Value = 804.0833
$ws.Range("I131").Value = 246.33333
$ws.Range("K131").Value = 738.99999
$ws.Range("M131").Value = 4301.00001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1100
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1100
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1100
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -4940
$ws.Range("H113").Value = 3452
$ws.Range("I113").Value = 3252.5
$ws.Range("K113").Value = 3252.5
$ws.Range("M113").Value = -1082.5
$ws.Range("H132").Value = 3445.3333
$ws.Range("I132").Value = 3932.5
$ws.Range("K132").Value = 11797.5
$ws.Range("M132").Value = -9267.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3499.5
$ws.Range("I40").Value = 3499.5
$ws.Range("K40").Value = 3499.5
$ws.Range("M40").Value = -3363.5
$ws.Range("H93").Value = 1352.3478
$ws.Range("I93").Value = 1244.5625
$ws.Range("J93").Value = 1598.7142
$ws.Range("K93").Value = 1244.5625
$ws.Range("L93").Value = 1598.7142
$ws.Range("M93").Value = 3.4375
$ws.Range("N93").Value = -4094.7142
$ws.Range("H133").Value = 49999.5
$ws.Range("J133").Value = 49999.5
$ws.Range("L133").Value = 49999.5
$ws.Range("N133").Value = -55059.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3998.1904
$ws.Range("I132").Value = 3325.2856
$ws.Range("K132").Value = 9975.856800000001
$ws.Range("M132").Value = -7445.856800000001
$ws.Range("H136").Value = 1517.0555
$ws.Range("I136").Value = 843.3570999999999
$ws.Range("J136").Value = 3875
$ws.Range("K136").Value = 2530.0713
$ws.Range("L136").Value = 11625
$ws.Range("M136").Value = 19.92870000000039
$ws.Range("N136").Value = -16725

Write-Host "Applied all Seraph_Profits updates"